$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 557, shifting existing rows 557-607 down to 558-608.
$ws.Rows.Item(557).Insert()

# Populate the newly inserted row 557 with the new record's data.
$ws.Cells.Item(557, 1).Value = 3
$ws.Cells.Item(557, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(557, 3).Value = "Coquimbo"
$ws.Cells.Item(557, 4).Value = 45132
$ws.Cells.Item(557, 5).Value = 5
$ws.Cells.Item(557, 6).Value = 100112043
$ws.Cells.Item(557, 7).Value = "Pepino ensalada"
$ws.Cells.Item(557, 8).Value = "Sin especificar"
$ws.Cells.Item(557, 9).Value = "Primera"
$ws.Cells.Item(557, 10).Value = 100
$ws.Cells.Item(557, 11).Value = 9500
$ws.Cells.Item(557, 12).Value = 10000
$ws.Cells.Item(557, 13).Value = 9750
$ws.Cells.Item(557, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(557, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(557, 16).Value = 162
$ws.Cells.Item(557, 17).Value = 60
$ws.Cells.Item(557, 18).Value = "Hortaliza"
